# Insert a new data row at row 358 in the "Tomate" subconjunto sheet.
# Excel shifts all rows from 358:439 down to 359:440, preserving their
# contents and formatting, leaving row 358 blank and ready to be filled
# with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("358:358").Insert()

$ws.Range("A358").Value = 5
$ws.Range("B358").Value = "Macroferia Regional de Talca"
$ws.Range("C358").Value = "Maule"
$ws.Range("D358").Value = 44543
$ws.Range("E358").Value = 7
$ws.Range("F358").Value = 100112020
$ws.Range("G358").Value = "Tomate"
$ws.Range("H358").Value = "Larga vida"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 3500
$ws.Range("K358").Value = 7000
$ws.Range("L358").Value = 7000
$ws.Range("M358").Value = 7000
$ws.Range("N358").Value = "$/caja 15 kilos"
$ws.Range("O358").Value = "Región del Maule"
$ws.Range("P358").Value = 467
$ws.Range("Q358").Value = 15
$ws.Range("R358").Value = "Hortaliza"
